# Deploying to gh-pages from @ codeforIATI/codelists@97b2f69c - SectorGroup.xlsx reorder
#
# The "codeforiati:*" metadata columns are re-ordered: the group-name/category-name
# columns move into the D/E position (swapping with the category-code/group-code
# columns, which move into the F/G position), for every row of the table, including
# the header row.
#
#   before:  D = category-code   E = group-code   F = group-name   G = category-name
#   after:   D = group-name      E = category-name  F = group-code   G = category-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$dgRng = $ws.Range("D1:G$lastRow")

# Ensure these columns keep being stored as text (the IATI sector/group codes such as
# "110" or "111" must not be silently converted into numbers).
$dgRng.NumberFormat = "@"

$arr = $dgRng.Value2

for ($r = 1; $r -le $lastRow; $r++) {
    $catCode   = $arr[$r,1]   # D (before): codeforiati:category-code
    $groupCode = $arr[$r,2]   # E (before): codeforiati:group-code
    $groupName = $arr[$r,3]   # F (before): codeforiati:group-name
    $catName   = $arr[$r,4]   # G (before): codeforiati:category-name

    $arr[$r,1] = $groupName   # D (after): codeforiati:group-name
    $arr[$r,2] = $catName     # E (after): codeforiati:category-name
    $arr[$r,3] = $groupCode   # F (after): codeforiati:group-code
    $arr[$r,4] = $catCode     # G (after): codeforiati:category-code
}

$dgRng.Value2 = $arr
